$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 ("b.md") moves from "Handed back: in sync with en-US"
#     to "Ready for handoff" in both the zh-cn and de-de status columns (B3, C3) ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: a new handoff was generated for b.md ---
$wsZh = $wb.Sheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-11 00:43:59"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$C`$3") {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: a new handoff was generated for b.md ---
$wsDe = $wb.Sheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-11 00:44:06"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$C`$3") {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
